# Auto-generated edit script: updates currentAveragePrice / leve price / profit
# columns (H-N) for specific leve rows across multiple worksheets, matching the
# scheduled-runner data refresh described in the commit message.
$wb = $excel.ActiveWorkbook

function Set-CellNumber($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Clear-Cell($ws, $row, $col) {
    $ws.Cells.Item($row, $col).ClearContents()
}

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
Set-CellNumber $ws 2 8 1044.8334  # H2
Set-CellNumber $ws 2 9 351.375  # I2
Set-CellNumber $ws 2 11 351.375  # K2
Set-CellNumber $ws 2 13 -238.375  # M2

# Row 12 (Leve Item ID 5515)
Set-CellNumber $ws 12 8 403  # H12
Set-CellNumber $ws 12 9 417.125  # I12
Set-CellNumber $ws 12 10 290  # J12
Set-CellNumber $ws 12 11 417.125  # K12
Set-CellNumber $ws 12 12 290  # L12
Set-CellNumber $ws 12 13 -247.125  # M12
Set-CellNumber $ws 12 14 -630  # N12

# Row 15 (Leve Item ID 44146)
Set-CellNumber $ws 15 8 1115.3478  # H15
Set-CellNumber $ws 15 9 1115.3478  # I15
Set-CellNumber $ws 15 11 3346.0434  # K15
Set-CellNumber $ws 15 13 -3177.0434  # M15

# Row 76 (Leve Item ID 12602)
Set-CellNumber $ws 76 8 6500  # H76
Set-CellNumber $ws 76 9 6000  # I76
Set-CellNumber $ws 76 11 6000  # K76
Set-CellNumber $ws 76 13 -5685  # M76

# Row 79 (Leve Item ID 12602)
Set-CellNumber $ws 79 8 6500  # H79
Set-CellNumber $ws 79 9 6000  # I79
Set-CellNumber $ws 79 11 6000  # K79
Set-CellNumber $ws 79 13 -4908  # M79

# Row 80 (Leve Item ID 12605)
Set-CellNumber $ws 80 8 537.63635  # H80
Set-CellNumber $ws 80 9 198.6  # I80
Set-CellNumber $ws 80 10 820.1667  # J80
Set-CellNumber $ws 80 11 595.8  # K80
Set-CellNumber $ws 80 12 2460.5001  # L80
Set-CellNumber $ws 80 13 402.2  # M80
Set-CellNumber $ws 80 14 -4456.5001  # N80

# Row 83 (Leve Item ID 12605)
Set-CellNumber $ws 83 8 537.63635  # H83
Set-CellNumber $ws 83 9 198.6  # I83
Set-CellNumber $ws 83 10 820.1667  # J83
Set-CellNumber $ws 83 11 1787.4  # K83
Set-CellNumber $ws 83 12 7381.5003  # L83
Set-CellNumber $ws 83 13 3204.6  # M83
Set-CellNumber $ws 83 14 -17365.5003  # N83

# Row 92 (Leve Item ID 19901)
Set-CellNumber $ws 92 8 891.625  # H92
Set-CellNumber $ws 92 9 838.8333  # I92
Set-CellNumber $ws 92 11 838.8333  # K92
Set-CellNumber $ws 92 13 409.1667  # M92

# Row 94 (Leve Item ID 19905)
Set-CellNumber $ws 94 8 541.6667  # H94
Set-CellNumber $ws 94 9 541.6667  # I94
Set-CellNumber $ws 94 11 541.6667  # K94
Set-CellNumber $ws 94 13 -90.66669999999999  # M94

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
Set-CellNumber $ws 32 8 5548.2095  # H32
Set-CellNumber $ws 32 9 5548.2095  # I32
Set-CellNumber $ws 32 10 0  # J32
Set-CellNumber $ws 32 11 5548.2095  # K32
Set-CellNumber $ws 32 12 0  # L32
Set-CellNumber $ws 32 13 -5261.2095  # M32
Clear-Cell $ws 32 14  # N32

# Row 46 (Leve Item ID 3498)
Set-CellNumber $ws 46 8 4746.5  # H46
Set-CellNumber $ws 46 9 5500  # I46
Set-CellNumber $ws 46 11 5500  # K46
Set-CellNumber $ws 46 13 -5181  # M46

# Row 55 (Leve Item ID 2830)
Set-CellNumber $ws 55 8 21453.572  # H55
Set-CellNumber $ws 55 10 29000  # J55
Set-CellNumber $ws 55 12 29000  # L55
Set-CellNumber $ws 55 14 -29630  # N55

# Row 61 (Leve Item ID 43999)
Set-CellNumber $ws 61 8 1690.7  # H61
Set-CellNumber $ws 61 9 1634.1111  # I61
Set-CellNumber $ws 61 10 2200  # J61
Set-CellNumber $ws 61 11 1634.1111  # K61
Set-CellNumber $ws 61 12 2200  # L61
Set-CellNumber $ws 61 13 -1422.1111  # M61
Set-CellNumber $ws 61 14 -2624  # N61

# Row 80 (Leve Item ID 10667)
Set-CellNumber $ws 80 8 38333.332  # H80

# Row 83 (Leve Item ID 10667)
Set-CellNumber $ws 83 8 38333.332  # H83

# Row 122 (Leve Item ID 36168)
Set-CellNumber $ws 122 8 2259.3125  # H122
Set-CellNumber $ws 122 9 2259.3125  # I122
Set-CellNumber $ws 122 11 6777.9375  # K122
Set-CellNumber $ws 122 13 -4327.9375  # M122

# Row 136 (Leve Item ID 43999)
Set-CellNumber $ws 136 8 1690.7  # H136
Set-CellNumber $ws 136 9 1634.1111  # I136
Set-CellNumber $ws 136 10 2200  # J136
Set-CellNumber $ws 136 11 4902.3333  # K136
Set-CellNumber $ws 136 12 6600  # L136
Set-CellNumber $ws 136 13 -2352.3333  # M136
Set-CellNumber $ws 136 14 -11700  # N136

$ws = $wb.Worksheets.Item("BSM")
# Row 82 (Leve Item ID 11877)
Set-CellNumber $ws 82 8 19607.533  # H82
Set-CellNumber $ws 82 9 6012.5557  # I82
Set-CellNumber $ws 82 11 6012.5557  # K82
Set-CellNumber $ws 82 13 -5629.5557  # M82

# Row 85 (Leve Item ID 11877)
Set-CellNumber $ws 85 8 19607.533  # H85
Set-CellNumber $ws 85 9 6012.5557  # I85
Set-CellNumber $ws 85 11 6012.5557  # K85
Set-CellNumber $ws 85 13 -4686.5557  # M85

# Row 86 (Leve Item ID 12526)
Set-CellNumber $ws 86 8 3640  # H86
Set-CellNumber $ws 86 9 3595.1667  # I86
Set-CellNumber $ws 86 11 3595.1667  # K86
Set-CellNumber $ws 86 13 -2472.1667  # M86

# Row 89 (Leve Item ID 12526)
Set-CellNumber $ws 89 8 3640  # H89
Set-CellNumber $ws 89 9 3595.1667  # I89
Set-CellNumber $ws 89 11 17975.8335  # K89
Set-CellNumber $ws 89 13 -12359.8335  # M89

# Row 105 (Leve Item ID 19947)
Set-CellNumber $ws 105 8 3090.3  # H105
Set-CellNumber $ws 105 9 2988  # I105
Set-CellNumber $ws 105 10 3243.75  # J105
Set-CellNumber $ws 105 11 2988  # K105
Set-CellNumber $ws 105 12 3243.75  # L105
Set-CellNumber $ws 105 13 -1241  # M105
Set-CellNumber $ws 105 14 -6737.75  # N105

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Leve Item ID 5367)
Set-CellNumber $ws 22 8 403.57144  # H22
Set-CellNumber $ws 22 9 106.25  # I22
Set-CellNumber $ws 22 11 106.25  # K22
Set-CellNumber $ws 22 13 243.75  # M22

# Row 41 (Leve Item ID 1917)
Set-CellNumber $ws 41 8 15500  # H41
Set-CellNumber $ws 41 10 21250  # J41
Set-CellNumber $ws 41 12 21250  # L41
Set-CellNumber $ws 41 14 -22106  # N41

# Row 59 (Leve Item ID 1942)
Set-CellNumber $ws 59 8 27786.363  # H59
Set-CellNumber $ws 59 9 21775  # I59
Set-CellNumber $ws 59 11 21775  # K59
Set-CellNumber $ws 59 13 -20630  # M59

# Row 60 (Leve Item ID 1937)
Set-CellNumber $ws 60 8 24996.875  # H60
Set-CellNumber $ws 60 10 24996.875  # J60
Set-CellNumber $ws 60 12 24996.875  # L60
Set-CellNumber $ws 60 14 -26018.875  # N60

# Row 107 (Leve Item ID 27689)
Set-CellNumber $ws 107 8 426.23077  # H107
Set-CellNumber $ws 107 9 261  # I107
Set-CellNumber $ws 107 11 261  # K107
Set-CellNumber $ws 107 13 1659  # M107

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
Set-CellNumber $ws 80 8 3967.0908  # H80
Set-CellNumber $ws 80 9 2299.3333  # I80
Set-CellNumber $ws 80 10 4592.5  # J80
Set-CellNumber $ws 80 11 2299.3333  # K80
Set-CellNumber $ws 80 12 4592.5  # L80
Set-CellNumber $ws 80 13 -1301.3333  # M80
Set-CellNumber $ws 80 14 -6588.5  # N80

# Row 83 (Leve Item ID 12521)
Set-CellNumber $ws 83 8 3967.0908  # H83
Set-CellNumber $ws 83 9 2299.3333  # I83
Set-CellNumber $ws 83 10 4592.5  # J83
Set-CellNumber $ws 83 11 11496.6665  # K83
Set-CellNumber $ws 83 12 22962.5  # L83
Set-CellNumber $ws 83 13 -6504.666499999999  # M83
Set-CellNumber $ws 83 14 -32946.5  # N83

# Row 97 (Leve Item ID 19940)
Set-CellNumber $ws 97 8 854.2105  # H97
Set-CellNumber $ws 97 9 890.05884  # I97
Set-CellNumber $ws 97 10 549.5  # J97
Set-CellNumber $ws 97 11 890.05884  # K97
Set-CellNumber $ws 97 12 549.5  # L97
Set-CellNumber $ws 97 13 -394.05884  # M97
Set-CellNumber $ws 97 14 -1541.5  # N97

# Row 102 (Leve Item ID 36169)
Set-CellNumber $ws 102 8 1448.8889  # H102
Set-CellNumber $ws 102 9 1445.8823  # I102
Set-CellNumber $ws 102 11 1445.8823  # K102
Set-CellNumber $ws 102 13 176.1177  # M102

# Row 113 (Leve Item ID 27710)
Set-CellNumber $ws 113 8 2299.6  # H113
Set-CellNumber $ws 113 9 1999.3334  # I113
Set-CellNumber $ws 113 10 2750  # J113
Set-CellNumber $ws 113 11 1999.3334  # K113
Set-CellNumber $ws 113 12 2750  # L113
Set-CellNumber $ws 113 13 170.6666  # M113
Set-CellNumber $ws 113 14 -7090  # N113

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
Set-CellNumber $ws 16 8 14342.667  # H16
Set-CellNumber $ws 16 10 38358  # J16
Set-CellNumber $ws 16 12 38358  # L16
Set-CellNumber $ws 16 14 -38698  # N16

# Row 22 (Leve Item ID 5277)
Set-CellNumber $ws 22 8 916  # H22
Set-CellNumber $ws 22 9 630  # I22
Set-CellNumber $ws 22 11 630  # K22
Set-CellNumber $ws 22 13 -335  # M22

# Row 27 (Leve Item ID 5277)
Set-CellNumber $ws 27 8 916  # H27
Set-CellNumber $ws 27 9 630  # I27
Set-CellNumber $ws 27 11 630  # K27
Set-CellNumber $ws 27 13 -523  # M27

# Row 46 (Leve Item ID 5282)
Set-CellNumber $ws 46 8 1858.8  # H46
Set-CellNumber $ws 46 9 900  # I46
Set-CellNumber $ws 46 10 3297  # J46
Set-CellNumber $ws 46 11 900  # K46
Set-CellNumber $ws 46 12 3297  # L46
Set-CellNumber $ws 46 13 -712  # M46
Set-CellNumber $ws 46 14 -3673  # N46

# Row 55 (Leve Item ID 5284)
Set-CellNumber $ws 55 8 1118.8889  # H55
Set-CellNumber $ws 55 9 1274  # I55
Set-CellNumber $ws 55 10 925  # J55
Set-CellNumber $ws 55 11 1274  # K55
Set-CellNumber $ws 55 12 925  # L55
Set-CellNumber $ws 55 13 -1101  # M55
Set-CellNumber $ws 55 14 -1271  # N55

# Row 132 (Leve Item ID 44058)
Set-CellNumber $ws 132 8 15747.45  # H132
Set-CellNumber $ws 132 9 15441.611  # I132
Set-CellNumber $ws 132 11 46324.833  # K132
Set-CellNumber $ws 132 13 -43794.833  # M132

$ws = $wb.Worksheets.Item("WVR")
# Row 15 (Leve Item ID 2670)
Set-CellNumber $ws 15 8 0  # H15
Set-CellNumber $ws 15 9 0  # I15
Set-CellNumber $ws 15 11 0  # K15
Clear-Cell $ws 15 13  # M15

# Row 43 (Leve Item ID 3831)
Set-CellNumber $ws 43 8 14000  # H43
Set-CellNumber $ws 43 10 14000  # J43
Set-CellNumber $ws 43 12 14000  # L43
Set-CellNumber $ws 43 14 -14298  # N43

# Row 81 (Leve Item ID 12596)
Set-CellNumber $ws 81 8 562  # H81
Set-CellNumber $ws 81 9 549.3333  # I81
Set-CellNumber $ws 81 10 600  # J81
Set-CellNumber $ws 81 11 1098.6666  # K81
Set-CellNumber $ws 81 12 1200  # L81
Set-CellNumber $ws 81 13 -37.66660000000002  # M81
Set-CellNumber $ws 81 14 -3322  # N81

# Row 84 (Leve Item ID 12596)
Set-CellNumber $ws 84 8 562  # H84
Set-CellNumber $ws 84 9 549.3333  # I84
Set-CellNumber $ws 84 10 600  # J84
Set-CellNumber $ws 84 11 5493.333000000001  # K84
Set-CellNumber $ws 84 12 6000  # L84
Set-CellNumber $ws 84 13 -189.3330000000005  # M84
Set-CellNumber $ws 84 14 -16608  # N84

# Row 107 (Leve Item ID 27746)
Set-CellNumber $ws 107 8 109.4  # H107
Set-CellNumber $ws 107 9 109.4  # I107
Set-CellNumber $ws 107 10 0  # J107
Set-CellNumber $ws 107 11 328.2  # K107
Set-CellNumber $ws 107 12 0  # L107
Set-CellNumber $ws 107 13 1591.8  # M107
Clear-Cell $ws 107 14  # N107

# Row 122 (Leve Item ID 36208)
Set-CellNumber $ws 122 8 1559  # H122
Set-CellNumber $ws 122 9 1559  # I122
Set-CellNumber $ws 122 11 4677  # K122
Set-CellNumber $ws 122 13 -2227  # M122

# Row 126 (Leve Item ID 36210)
Set-CellNumber $ws 126 8 1715.1666  # H126
Set-CellNumber $ws 126 9 1878.2  # I126
Set-CellNumber $ws 126 10 900  # J126
Set-CellNumber $ws 126 11 5634.6  # K126
Set-CellNumber $ws 126 12 2700  # L126
Set-CellNumber $ws 126 13 -3164.6  # M126
Set-CellNumber $ws 126 14 -7640  # N126

# Row 132 (Leve Item ID 44029)
Set-CellNumber $ws 132 8 738  # H132
Set-CellNumber $ws 132 9 757  # I132
Set-CellNumber $ws 132 10 700  # J132
Set-CellNumber $ws 132 11 2271  # K132
Set-CellNumber $ws 132 12 2100  # L132
Set-CellNumber $ws 132 13 259  # M132
Set-CellNumber $ws 132 14 -7160  # N132
